$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "K" (column G) values regenerated from source (Strike# -> K), row => value
$kValues = @{
    2  = 1
    3  = 2
    4  = 0
    5  = 2
    6  = 1
    7  = 2
    8  = 1
    9  = 1
    10 = 0
    11 = 0
    12 = 1
    13 = 1
    14 = 1
    15 = 2
    16 = 0
    17 = 1
    18 = 1
    19 = 3
    20 = 2
    21 = 2
    22 = 0
    23 = 1
    24 = 0
    25 = 1
    26 = 1
    27 = 0
    28 = 0
    29 = 1
    30 = 1
    31 = 3
    32 = 2
    33 = 0
    34 = 0
    36 = 2
    37 = 1
    38 = 0
    39 = 2
    40 = 1
    41 = 1
    42 = 2
    43 = 1
    45 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
